$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("Rules" sheet) holds the rule-name for the 4th time band.
# It was "R40"; change it to the text "1" (still a plain text cell,
# not a number) while leaving the cell's existing formatting/style
# untouched.
#
# A direct  $ws.Range("B11").Value = "1"  would be auto-coerced to the
# *number* 1 by Excel, and prefixing with an apostrophe to force text
# stamps a quote-prefix flag onto the cell's style. Instead, enter a
# formula that evaluates to the text "1", then convert it in place to
# a plain value via copy / paste-special-values - this yields a true
# text constant with no formula and no style change.
$target = $ws.Range("B11")
$target.Formula = '="1"'
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
